$d = $word.ActiveDocument

# Locate the run containing " se unosi i sistem i " (the text is entirely
# inside one run, so the Find range exactly bounds that run).
$hit = $d.Content
$found = $hit.Find.Execute(" se unosi i sistem i ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target text ' se unosi i sistem i '"
}
$start = $hit.Start
$end = $hit.End

# Fix "i" (and) -> "u" (into): "... se unosi i sistem i ..." -> "... se unosi u sistem i ...".
# (This single assignment also normalises/merges the runs that follow in
# the paragraph, which we re-split back to their original boundaries below.)
$whole = $d.Range($start, $end)
$whole.Text = " se unosi u sistem i "

# Character offsets (relative to $start / $end) for every run boundary that
# needs to exist afterwards:
#   " se unosi " | "u" | " sistem i " | "njegova lozinka" | " je promenjen" | "a"
$p1 = $start + 10   # end of " se unosi "
$p2 = $start + 11   # end of "u"
$p3 = $end + 15     # end of "njegova lozinka"
$p4 = $end + 28     # end of " je promenjen"
$p5 = $end + 29     # end of "a"

# Touching (and reverting) a character formatting property on each
# sub-range forces the engine to keep/create a distinct <w:r> run at that
# boundary instead of merging adjacent same-formatted text back together.
$s1 = $d.Range($start, $p1)
$s1.Bold = 1
$s1.Bold = 0

$s2 = $d.Range($p1, $p2)
$s2.Bold = 1
$s2.Bold = 0

$s3 = $d.Range($p2, $end)
$s3.Bold = 1
$s3.Bold = 0

$s4 = $d.Range($end, $p3)
$s4.Bold = 1
$s4.Bold = 0

$s5 = $d.Range($p3, $p4)
$s5.Bold = 1
$s5.Bold = 0

$s6 = $d.Range($p4, $p5)
$s6.Bold = 1
$s6.Bold = 0
